# Clean up test data: remove/update paths pointing to /node/xxx,
# and refresh the saved view state (active sheet/selection) to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("pages_with_related_resources")
$ws3 = $wb.Worksheets.Item("pages_with_external_resources")

# Fix the /node/ path that appears on both the
# "pages_with_related_resources" and "pages_with_external_resources" sheets.
$newPath = "espanol/news-events/cancer-currents-blog/2019/vitamina-d-complemento-cancer-prevencion"
$ws1.Range("A4").Value = $newPath
$ws3.Range("A9").Value = $newPath

# Update saved selection/active sheet state.
$ws1.Activate()
$ws1.Range("A4").Select()

$ws3.Activate()
$ws3.Range("A16").Select()
